$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2982.4167
$ws.Range("I64").Value = 3004.7778
$ws.Range("J64").Value = 2915.3333
$ws.Range("K64").Value = 3004.7778
$ws.Range("L64").Value = 2915.3333
$ws.Range("M64").Value = -2756.7778
$ws.Range("N64").Value = -3411.3333

$ws.Range("H67").Value = 2982.4167
$ws.Range("I67").Value = 3004.7778
$ws.Range("J67").Value = 2915.3333
$ws.Range("K67").Value = 3004.7778
$ws.Range("L67").Value = 2915.3333
$ws.Range("M67").Value = -2146.7778
$ws.Range("N67").Value = -4631.3333

$ws.Range("H125").Value = 3598
$ws.Range("I125").Value = 1894
$ws.Range("J125").Value = 4450
$ws.Range("K125").Value = 17046
$ws.Range("L125").Value = 40050
$ws.Range("M125").Value = -14586
$ws.Range("N125").Value = -44970

$ws.Range("H129").Value = 4099358
$ws.Range("J129").Value = 956.1607
$ws.Range("L129").Value = 2868.4821
$ws.Range("N129").Value = -12868.4821

$ws.Range("H137").Value = 2259.5
$ws.Range("I137").Value = 2251.6365
$ws.Range("K137").Value = 6754.9095
$ws.Range("M137").Value = -4204.9095

$ws.Range("H141").Value = 349873.38
$ws.Range("I141").Value = 1298.6666
$ws.Range("J141").Value = 1918459.5
$ws.Range("K141").Value = 3895.9998
$ws.Range("L141").Value = 5755378.5
$ws.Range("M141").Value = 1284.0002
$ws.Range("N141").Value = -5765738.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 994.0909
$ws.Range("I74").Value = 905
$ws.Range("J74").Value = 1150
$ws.Range("K74").Value = 905
$ws.Range("L74").Value = 1150
$ws.Range("M74").Value = -31
$ws.Range("N74").Value = -2898

$ws.Range("H77").Value = 994.0909
$ws.Range("I77").Value = 905
$ws.Range("J77").Value = 1150
$ws.Range("K77").Value = 4525
$ws.Range("L77").Value = 5750
$ws.Range("M77").Value = -157
$ws.Range("N77").Value = -14486

$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2266.923
$ws.Range("I99").Value = 584.3333
$ws.Range("J99").Value = 6052.75
$ws.Range("K99").Value = 584.3333
$ws.Range("L99").Value = 6052.75
$ws.Range("M99").Value = 913.6667
$ws.Range("N99").Value = -9048.75

$ws.Range("H107").Value = 4042.6
$ws.Range("I107").Value = 2459
$ws.Range("J107").Value = 7737.6665
$ws.Range("K107").Value = 2459
$ws.Range("L107").Value = 7737.6665
$ws.Range("M107").Value = -539
$ws.Range("N107").Value = -11577.6665

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2649.608
$ws.Range("I31").Value = 1641.44
$ws.Range("J31").Value = 3619
$ws.Range("K31").Value = 1641.44
$ws.Range("L31").Value = 3619
$ws.Range("M31").Value = -1346.44
$ws.Range("N31").Value = -4209

$ws.Range("H34").Value = 2649.608
$ws.Range("I34").Value = 1641.44
$ws.Range("J34").Value = 3619
$ws.Range("K34").Value = 1641.44
$ws.Range("L34").Value = 3619
$ws.Range("M34").Value = -1439.44
$ws.Range("N34").Value = -4023

$ws.Range("H63").Value = 39600
$ws.Range("J63").Value = 39600
$ws.Range("L63").Value = 39600
$ws.Range("N63").Value = -40972

$ws.Range("H66").Value = 39600
$ws.Range("J66").Value = 39600
$ws.Range("L66").Value = 118800
$ws.Range("N66").Value = -125664

$ws.Range("H74").Value = 18541.111
$ws.Range("J74").Value = 18541.111
$ws.Range("L74").Value = 18541.111
$ws.Range("N74").Value = -20289.111

$ws.Range("H77").Value = 18541.111
$ws.Range("J77").Value = 18541.111
$ws.Range("L77").Value = 55623.333
$ws.Range("N77").Value = -64359.333

$ws.Range("H134").Value = 2612.6843
$ws.Range("I134").Value = 1070.1111
$ws.Range("J134").Value = 4001
$ws.Range("K134").Value = 3210.3333
$ws.Range("L134").Value = 12003
$ws.Range("M134").Value = -675.3333000000002
$ws.Range("N134").Value = -17073

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3975.889
$ws.Range("I68").Value = 450
$ws.Range("J68").Value = 6796.6
$ws.Range("K68").Value = 1350
$ws.Range("L68").Value = 20389.8
$ws.Range("M68").Value = -539
$ws.Range("N68").Value = -22011.8

$ws.Range("H71").Value = 3975.889
$ws.Range("I71").Value = 450
$ws.Range("J71").Value = 6796.6
$ws.Range("K71").Value = 4050
$ws.Range("L71").Value = 61169.4
$ws.Range("M71").Value = 6
$ws.Range("N71").Value = -69281.39999999999

$ws.Range("H110").Value = 3300.7778
$ws.Range("I110").Value = 1909
$ws.Range("J110").Value = 3996.6667
$ws.Range("K110").Value = 5727
$ws.Range("L110").Value = 11990.0001
$ws.Range("M110").Value = -1637
$ws.Range("N110").Value = -20170.0001

$ws.Range("H117").Value = 1804.5714
$ws.Range("I117").Value = 300
$ws.Range("J117").Value = 2933
$ws.Range("K117").Value = 900
$ws.Range("L117").Value = 8799
$ws.Range("M117").Value = 2542
$ws.Range("N117").Value = -15683

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1907.5385
$ws.Range("J46").Value = 2822.5
$ws.Range("L46").Value = 2822.5
$ws.Range("N46").Value = -3198.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 34353.5
$ws.Range("J46").Value = 34353.5
$ws.Range("L46").Value = 34353.5
$ws.Range("N46").Value = -34815.5

$ws.Range("H132").Value = 7343.625
$ws.Range("I132").Value = 1468.2821
$ws.Range("K132").Value = 4404.846299999999
$ws.Range("M132").Value = -1874.846299999999

$ws.Range("H134").Value = 34353.5
$ws.Range("J134").Value = 34353.5
$ws.Range("L134").Value = 103060.5
$ws.Range("N134").Value = -108130.5

$ws.Range("H140").Value = 62085.4
$ws.Range("J140").Value = 62085.4
$ws.Range("L140").Value = 62085.4
$ws.Range("N140").Value = -72445.39999999999
